# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows for "Black Amber" (Especial/Primera/Segunda)
# right before the existing row 29, shifting the remaining data down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at position 29 (existing rows 29:110 shift down to 32:113)
$ws.Rows("29:31").Insert()

# Common (constant) column values used across the whole data table
$mercadoId = 9
$mercado   = "Vega Central Mapocho de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$tipo      = "Fruta"
$prodId    = 100103
$prodNom   = "Frutos de hueso (carozo)"
$catId     = 100103002
$categoria = "Ciruela"

# Row 29: Black Amber - Especial
$ws.Range("A29").Value = $mercadoId
$ws.Range("B29").Value = $mercado
$ws.Range("C29").Value = $region
$ws.Range("D29").Value = 44624
$ws.Range("E29").Value = $codreg
$ws.Range("F29").Value = $tipo
$ws.Range("G29").Value = $prodId
$ws.Range("H29").Value = $prodNom
$ws.Range("I29").Value = $catId
$ws.Range("J29").Value = $categoria
$ws.Range("K29").Value = "Black Amber"
$ws.Range("L29").Value = "Especial"
$ws.Range("M29").Value = 210
$ws.Range("N29").Value = 11200
$ws.Range("O29").Value = 11200
$ws.Range("P29").Value = 11200
$ws.Range("Q29").Value = "$/caja 16 kilos granel"
$ws.Range("R29").Value = "Región de O'Higgins"
$ws.Range("S29").Value = 700
$ws.Range("T29").Value = 16

# Row 30: Black Amber - Primera
$ws.Range("A30").Value = $mercadoId
$ws.Range("B30").Value = $mercado
$ws.Range("C30").Value = $region
$ws.Range("D30").Value = 44624
$ws.Range("E30").Value = $codreg
$ws.Range("F30").Value = $tipo
$ws.Range("G30").Value = $prodId
$ws.Range("H30").Value = $prodNom
$ws.Range("I30").Value = $catId
$ws.Range("J30").Value = $categoria
$ws.Range("K30").Value = "Black Amber"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 250
$ws.Range("N30").Value = 9600
$ws.Range("O30").Value = 9600
$ws.Range("P30").Value = 9600
$ws.Range("Q30").Value = "$/caja 16 kilos granel"
$ws.Range("R30").Value = "Región de O'Higgins"
$ws.Range("S30").Value = 600
$ws.Range("T30").Value = 16

# Row 31: Black Amber - Segunda
$ws.Range("A31").Value = $mercadoId
$ws.Range("B31").Value = $mercado
$ws.Range("C31").Value = $region
$ws.Range("D31").Value = 44624
$ws.Range("E31").Value = $codreg
$ws.Range("F31").Value = $tipo
$ws.Range("G31").Value = $prodId
$ws.Range("H31").Value = $prodNom
$ws.Range("I31").Value = $catId
$ws.Range("J31").Value = $categoria
$ws.Range("K31").Value = "Black Amber"
$ws.Range("L31").Value = "Segunda"
$ws.Range("M31").Value = 280
$ws.Range("N31").Value = 8000
$ws.Range("O31").Value = 8000
$ws.Range("P31").Value = 8000
$ws.Range("Q31").Value = "$/caja 16 kilos granel"
$ws.Range("R31").Value = "Región de O'Higgins"
$ws.Range("S31").Value = 500
$ws.Range("T31").Value = 16
